$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns at D:E, shifting old D:K data to F:M
$ws.Columns("D:E").Insert()

# Step 2: Copy formatting (number format, font, style) from F:M (shifted original columns)
# into the newly inserted D:E columns so they match (date style for row7/38/80, number style otherwise)
$ws.Range("F8:G102").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("F7:G7").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("F38:G38").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("F80:G80").Copy()
$ws.Range("D80").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Step 3: Write new quarter values (columns D and E) and correct a handful of restated historical values
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("D8").Value = 406400
$ws.Range("E8").Value = 436000
$ws.Range("D9").Value = 68600
$ws.Range("E9").Value = 64900
$ws.Range("D10").Value = 337800
$ws.Range("E10").Value = 371100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 1800
$ws.Range("E15").Value = 2200
$ws.Range("D17").Value = 285300
$ws.Range("E17").Value = 291500
$ws.Range("D18").Value = 121100
$ws.Range("E18").Value = 144500
$ws.Range("D20").Value = 11300
$ws.Range("E20").Value = 12600
$ws.Range("D21").Value = 139000
$ws.Range("E21").Value = 163700
$ws.Range("D22").Value = 14500
$ws.Range("E22").Value = 17600
$ws.Range("D23").Value = 117900
$ws.Range("E23").Value = 139500
$ws.Range("D24").Value = 27600
$ws.Range("E24").Value = 37500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 90300
$ws.Range("E26").Value = 102000
$ws.Range("D27").Value = 86800
$ws.Range("E27").Value = 104800
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 700
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -11300
$ws.Range("E32").Value = -12600
$ws.Range("D33").Value = 86800
$ws.Range("E33").Value = 105500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 86800
$ws.Range("E35").Value = 105500
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("D41").Value = 495100
$ws.Range("E41").Value = 817300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 223900
$ws.Range("E43").Value = 236700
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 2061700
$ws.Range("E47").Value = 1957900
$ws.Range("D48").Value = 60100
$ws.Range("E48").Value = 52400
$ws.Range("D49").Value = 338700
$ws.Range("E49").Value = 340600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 91000
$ws.Range("E52").Value = 94500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3336100
$ws.Range("E54").Value = 3599300
$ws.Range("D57").Value = 80000
$ws.Range("E57").Value = 91400
$ws.Range("D58").Value = 68500
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 126200
$ws.Range("E59").Value = 285600
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 1460800
$ws.Range("E61").Value = 1492700
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 2268400
$ws.Range("E66").Value = 2491900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 1131100
$ws.Range("E72").Value = 1150700
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1067700
$ws.Range("E76").Value = 1107400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("D81").Value = 86800
$ws.Range("E81").Value = 105500
$ws.Range("D83").Value = 6600
$ws.Range("E83").Value = 6600
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 34300
$ws.Range("E89").Value = -17600
$ws.Range("D91").Value = -8300
$ws.Range("E91").Value = -5900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -283100
$ws.Range("E94").Value = -172500
$ws.Range("D96").Value = -43200
$ws.Range("E96").Value = -35300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -75000
$ws.Range("E100").Value = 229300
$ws.Range("D101").Value = 1100
$ws.Range("E101").Value = -1400
$ws.Range("D102").Value = -322600
$ws.Range("E102").Value = 37800
$ws.Range("H8").Value = 420200
$ws.Range("H9").Value = 73000
$ws.Range("H10").Value = 347200
$ws.Range("H17").Value = 284400
$ws.Range("H20").Value = 4400
$ws.Range("H32").Value = -4400
$ws.Range("G89").Value = 66800
$ws.Range("H89").Value = -59200
$ws.Range("I89").Value = -19200
$ws.Range("I91").Value = -4400
$ws.Range("J91").Value = -2300
$ws.Range("G102").Value = -22300
$ws.Range("H102").Value = -76500
